# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Update DAMSLTag (column I) and DialogAct (column J) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 4;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 12; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 13; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 16; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 21; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 22; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 24; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 31; Tag = "qy"; Act = "Yes-No-Question" },
    @{ Row = 42; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 47; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 67; Tag = "b";  Act = "Acknowledge (Backchannel)" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
